$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Highlight (yellow fill) existing cells that were flagged ---
# Column C, rows 2-7
$ws.Range("C2:C7").Interior.Color = 65535

# Column A, rows 8-13 and 15-17 (row 14 stays unhighlighted)
$ws.Range("A8:A13").Interior.Color = 65535
$ws.Range("A15:A17").Interior.Color = 65535

# --- Append 14 new deployment rows (18-31) ---
# First, give the new rows the same base formatting as the existing data rows
# (font styling matches the rest of the table) by copying formats down.
$ws.Range("A2").Copy()
$ws.Range("A18:A31").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B18:B31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ids = @("A004_SD013","A008_SD021","A009_SD012","A011_SD018","A014_SD005","A021_SD003","A024_SD004","A001_SD007","A002_SD030","A007_SD008","A010_SD019","A017_SD017","A025_SD006","A025_SD100")

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = 18 + $i
    $ws.Cells.Item($row, 1).Value = $ids[$i]
    $ws.Cells.Item($row, 2).Value = "Not logged"
}

# Highlight the new identifiers in column A just like the other flagged rows
$ws.Range("A18:A31").Interior.Color = 65535

# --- Restore the active selection like the source file ---
$ws.Range("A17").Select() | Out-Null
